# fix: 1/12/2025: Fix the algorithm/conditions on filtering the status of candidates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 5, pushing the existing rows 5-9 down to 8-12.
# Inserting above row 5 copies formatting from row 5 (the row being pushed down),
# which matches the desired formatting (style "2" with date format) for column F.
$ws.Rows("5:7").Insert()

# New row 5: Dash0 / Sean Guillen
$ws.Cells.Item(5, 1).Value = 708
$ws.Cells.Item(5, 2).Value = "Dash0"
$ws.Cells.Item(5, 3).Value = "Sales Engineer (US) x 3"
$ws.Cells.Item(5, 4).Value = "Sean Guillen"
$ws.Cells.Item(5, 5).Value = "2nd Interview"
$ws.Cells.Item(5, 6).Value = 45973

# New row 6: Rox / Alex Biller
$ws.Cells.Item(6, 1).Value = 727
$ws.Cells.Item(6, 2).Value = "Rox"
$ws.Cells.Item(6, 3).Value = "RVP Sales West (SF)"
$ws.Cells.Item(6, 4).Value = "Alex Biller"
$ws.Cells.Item(6, 5).Value = "2nd Interview"
$ws.Cells.Item(6, 6).Value = 45989

# New row 7: PointFive / Tony Bermeo
$ws.Cells.Item(7, 1).Value = 730
$ws.Cells.Item(7, 2).Value = "PointFive"
$ws.Cells.Item(7, 3).Value = "PointFive SE EST"
$ws.Cells.Item(7, 4).Value = "Tony Bermeo"
$ws.Cells.Item(7, 5).Value = "1st Interview"
$ws.Cells.Item(7, 6).Value = 45967
